# Apply updated "想去人数" (F column) counts per commit 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 64
$ws.Range("F6").Value = 3820
$ws.Range("F8").Value = 2538
$ws.Range("F9").Value = 73
$ws.Range("F10").Value = 3089
$ws.Range("F12").Value = 531
$ws.Range("F13").Value = 2295
$ws.Range("F15").Value = 115
$ws.Range("F17").Value = 446
$ws.Range("F22").Value = 305
$ws.Range("F23").Value = 363
$ws.Range("F24").Value = 653
$ws.Range("F26").Value = 43
$ws.Range("F29").Value = 147
$ws.Range("F32").Value = 43
$ws.Range("F33").Value = 4256
$ws.Range("F34").Value = 3941
$ws.Range("F36").Value = 5
$ws.Range("F38").Value = 1118
$ws.Range("F45").Value = 96
$ws.Range("F48").Value = 57

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2272

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value = 64
$ws.Range("F10").Value = 3820
$ws.Range("F12").Value = 2538
$ws.Range("F13").Value = 73
$ws.Range("F14").Value = 3089
$ws.Range("F15").Value = 531
$ws.Range("F16").Value = 2295
$ws.Range("F18").Value = 115
$ws.Range("F20").Value = 446
$ws.Range("F23").Value = 305
$ws.Range("F24").Value = 363
$ws.Range("F25").Value = 653
$ws.Range("F27").Value = 43
$ws.Range("F30").Value = 147
$ws.Range("F32").Value = 43
$ws.Range("F35").Value = 4256
$ws.Range("F36").Value = 3941
$ws.Range("F38").Value = 1118
$ws.Range("F46").Value = 96
$ws.Range("F48").Value = 57
